$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,7

$data[0,0] = 1
$data[0,1] = "BTC"
$data[0,2] = "Bitcoin"
$data[0,3] = 34312
$data[0,4] = 670427960957
$data[0,5] = 7289511810
$data[0,6] = 0.52088
$data[1,0] = 2
$data[1,1] = "ETH"
$data[1,2] = "Ethereum"
$data[1,3] = 1790
$data[1,4] = 215363841937
$data[1,5] = 5305371011
$data[1,6] = 0.29337
$data[2,0] = 3
$data[2,1] = "USDT"
$data[2,2] = "Tether"
$data[2,3] = 1
$data[2,4] = 84618317471
$data[2,5] = 13296117472
$data[2,6] = -0.02686
$data[3,0] = 4
$data[3,1] = "BNB"
$data[3,2] = "BNB"
$data[3,3] = 226.19
$data[3,4] = 34800995032
$data[3,5] = 266714652
$data[3,6] = 0.09859999999999999
$data[4,0] = 5
$data[4,1] = "XRP"
$data[4,2] = "XRP"
$data[4,3] = 0.556534
$data[4,4] = 29786263009
$data[4,5] = 567077708
$data[4,6] = 1.88312
$data[5,0] = 6
$data[5,1] = "USDC"
$data[5,2] = "USDC"
$data[5,3] = 1
$data[5,4] = 24989383116
$data[5,5] = 2972893263
$data[5,6] = 0.07561
$data[6,0] = 7
$data[6,1] = "STETH"
$data[6,2] = "Lido Staked Ether"
$data[6,3] = 1789.78
$data[6,4] = 15732806790
$data[6,5] = 4009453
$data[6,6] = 0.262
$data[7,0] = 8
$data[7,1] = "SOL"
$data[7,2] = "Solana"
$data[7,3] = 32.23
$data[7,4] = 13511907624
$data[7,5] = 461434586
$data[7,6] = 0.94301
$data[8,0] = 9
$data[8,1] = "ADA"
$data[8,2] = "Cardano"
$data[8,3] = 0.295553
$data[8,4] = 10330747553
$data[8,5] = 121520297
$data[8,6] = 1.32929
$data[9,0] = 10
$data[9,1] = "DOGE"
$data[9,2] = "Dogecoin"
$data[9,3] = 0.068906
$data[9,4] = 9762582123
$data[9,5] = 254676570
$data[9,6] = -0.03868
$data[10,0] = 11
$data[10,1] = "TRX"
$data[10,2] = "TRON"
$data[10,3] = 0.09482599999999999
$data[10,4] = 8429182348
$data[10,5] = 180000698
$data[10,6] = 0.5876
$data[11,0] = 12
$data[11,1] = "TON"
$data[11,2] = "Toncoin"
$data[11,3] = 2.06
$data[11,4] = 7124038571
$data[11,5] = 5748835
$data[11,6] = 0.88487
$data[12,0] = 13
$data[12,1] = "LINK"
$data[12,2] = "Chainlink"
$data[12,3] = 10.96
$data[12,4] = 6093250598
$data[12,5] = 452217852
$data[12,6] = -2.31299
$data[13,0] = 14
$data[13,1] = "MATIC"
$data[13,2] = "Polygon"
$data[13,3] = 0.626914
$data[13,4] = 5799087071
$data[13,5] = 127574608
$data[13,6] = 1.23815
$data[14,0] = 15
$data[14,1] = "WBTC"
$data[14,2] = "Wrapped Bitcoin"
$data[14,3] = 34304
$data[14,4] = 5622338627
$data[14,5] = 70012413
$data[14,6] = 0.50054
$data[15,0] = 16
$data[15,1] = "DOT"
$data[15,2] = "Polkadot"
$data[15,3] = 4.2
$data[15,4] = 5419861274
$data[15,5] = 107760588
$data[15,6] = 0.15275
$data[16,0] = 17
$data[16,1] = "LTC"
$data[16,2] = "Litecoin"
$data[16,3] = 68.06
$data[16,4] = 5026013898
$data[16,5] = 272887019
$data[16,6] = 0.19949
$data[17,0] = 18
$data[17,1] = "BCH"
$data[17,2] = "Bitcoin Cash"
$data[17,3] = 246.37
$data[17,4] = 4821140562
$data[17,5] = 118671321
$data[17,6] = 0.11192
$data[18,0] = 19
$data[18,1] = "SHIB"
$data[18,2] = "Shiba Inu"
$data[18,3] = 0.00000801
$data[18,4] = 4721249722
$data[18,5] = 158580199
$data[18,6] = 2.69355
$data[19,0] = 20
$data[19,1] = "AVAX"
$data[19,2] = "Avalanche"
$data[19,3] = 10.94
$data[19,4] = 3882843494
$data[19,5] = 145213718
$data[19,6] = 1.05547
$data[20,0] = 21
$data[20,1] = "DAI"
$data[20,2] = "Dai"
$data[20,3] = 0.998963
$data[20,4] = 3742177861
$data[20,5] = 73001340
$data[20,6] = -0.07834000000000001
$data[21,0] = 22
$data[21,1] = "LEO"
$data[21,2] = "LEO Token"
$data[21,3] = 4
$data[21,4] = 3713293536
$data[21,5] = 255942
$data[21,6] = 1.02362
$data[22,0] = 23
$data[22,1] = "TUSD"
$data[22,2] = "TrueUSD"
$data[22,3] = 1
$data[22,4] = 3361113405
$data[22,5] = 187773264
$data[22,6] = -0.00718
$data[23,0] = 24
$data[23,1] = "XLM"
$data[23,2] = "Stellar"
$data[23,3] = 0.115538
$data[23,4] = 3215715246
$data[23,5] = 44207029
$data[23,6] = 1.58621
$data[24,0] = 25
$data[24,1] = "UNI"
$data[24,2] = "Uniswap"
$data[24,3] = 4.17
$data[24,4] = 3138489478
$data[24,5] = 104711873
$data[24,6] = 1.56608
$data[25,0] = 26
$data[25,1] = "XMR"
$data[25,2] = "Monero"
$data[25,3] = 162.21
$data[25,4] = 2944524217
$data[25,5] = 57844779
$data[25,6] = 0.66301
$data[26,0] = 27
$data[26,1] = "OKB"
$data[26,2] = "OKB"
$data[26,3] = 44.87
$data[26,4] = 2692308748
$data[26,5] = 2594793
$data[26,6] = -0.25991
$data[27,0] = 28
$data[27,1] = "ETC"
$data[27,2] = "Ethereum Classic"
$data[27,3] = 16.35
$data[27,4] = 2342375469
$data[27,5] = 86432836
$data[27,6] = 0.35726
$data[28,0] = 29
$data[28,1] = "ATOM"
$data[28,2] = "Cosmos Hub"
$data[28,3] = 7.17
$data[28,4] = 2098086872
$data[28,5] = 77202991
$data[28,6] = 0.33921
$data[29,0] = 30
$data[29,1] = "BUSD"
$data[29,2] = "BUSD"
$data[29,3] = 1.001
$data[29,4] = 2057461892
$data[29,5] = 1089222194
$data[29,6] = -0.02034
$data[30,0] = 31
$data[30,1] = "HBAR"
$data[30,2] = "Hedera"
$data[30,3] = 0.05198
$data[30,4] = 1743188688
$data[30,5] = 28586125
$data[30,6] = 0.18661
$data[31,0] = 32
$data[31,1] = "FIL"
$data[31,2] = "Filecoin"
$data[31,3] = 3.77
$data[31,4] = 1742400804
$data[31,5] = 81970210
$data[31,6] = 2.91673
$data[32,0] = 33
$data[32,1] = "ICP"
$data[32,2] = "Internet Computer"
$data[32,3] = 3.9
$data[32,4] = 1737586813
$data[32,5] = 36349825
$data[32,6] = 7.81084
$data[33,0] = 34
$data[33,1] = "APT"
$data[33,2] = "Aptos"
$data[33,3] = 6.68
$data[33,4] = 1651969900
$data[33,5] = 83649689
$data[33,6] = 0.19053
$data[34,0] = 35
$data[34,1] = "LDO"
$data[34,2] = "Lido DAO"
$data[34,3] = 1.81
$data[34,4] = 1608679194
$data[34,5] = 17415570
$data[34,6] = -0.24381
$data[35,0] = 36
$data[35,1] = "CRO"
$data[35,2] = "Cronos"
$data[35,3] = 0.058995
$data[35,4] = 1553667644
$data[35,5] = 7224861
$data[35,6] = -0.81924
$data[36,0] = 37
$data[36,1] = "QNT"
$data[36,2] = "Quant"
$data[36,3] = 105.41
$data[36,4] = 1532154903
$data[36,5] = 16603066
$data[36,6] = -2.05444
$data[37,0] = 38
$data[37,1] = "VET"
$data[37,2] = "VeChain"
$data[37,3] = 0.01900394
$data[37,4] = 1381922766
$data[37,5] = 31984549
$data[37,6] = -1.26509
$data[38,0] = 39
$data[38,1] = "MKR"
$data[38,2] = "Maker"
$data[38,3] = 1436.92
$data[38,4] = 1296806385
$data[38,5] = 68243485
$data[38,6] = -0.76574
$data[39,0] = 40
$data[39,1] = "OP"
$data[39,2] = "Optimism"
$data[39,3] = 1.38
$data[39,4] = 1212599841
$data[39,5] = 47806352
$data[39,6] = 0.06126
$data[40,0] = 41
$data[40,1] = "NEAR"
$data[40,2] = "NEAR Protocol"
$data[40,3] = 1.22
$data[40,4] = 1209576589
$data[40,5] = 52886797
$data[40,6] = 0.23659
$data[41,0] = 42
$data[41,1] = "AAVE"
$data[41,2] = "Aave"
$data[41,3] = 82.09999999999999
$data[41,4] = 1198453475
$data[41,5] = 77616671
$data[41,6] = 2.73772
$data[42,0] = 43
$data[42,1] = "ARB"
$data[42,2] = "Arbitrum"
$data[42,3] = 0.926768
$data[42,4] = 1181616421
$data[42,5] = 107877458
$data[42,6] = 0.82775
$data[43,0] = 44
$data[43,1] = "MNT"
$data[43,2] = "Mantle"
$data[43,3] = 0.379514
$data[43,4] = 1179222082
$data[43,5] = 24218652
$data[43,6] = 0.13123
$data[44,0] = 45
$data[44,1] = "INJ"
$data[44,2] = "Injective"
$data[44,3] = 14.05
$data[44,4] = 1178185684
$data[44,5] = 77864110
$data[44,6] = 5.11589
$data[45,0] = 46
$data[45,1] = "KAS"
$data[45,2] = "Kaspa"
$data[45,3] = 0.051999
$data[45,4] = 1108233448
$data[45,5] = 8548916
$data[45,6] = 1.9511
$data[46,0] = 47
$data[46,1] = "RETH"
$data[46,2] = "Rocket Pool ETH"
$data[46,3] = 1946.53
$data[46,4] = 1042370166
$data[46,5] = 5309743
$data[46,6] = 0.04407
$data[47,0] = 48
$data[47,1] = "BSV"
$data[47,2] = "Bitcoin SV"
$data[47,3] = 50.62
$data[47,4] = 991091310
$data[47,5] = 68993912
$data[47,6] = 3.44628
$data[48,0] = 49
$data[48,1] = "RNDR"
$data[48,2] = "Render"
$data[48,3] = 2.61
$data[48,4] = 974275092
$data[48,5] = 102135896
$data[48,6] = 8.50577
$data[49,0] = 50
$data[49,1] = "GRT"
$data[49,2] = "The Graph"
$data[49,3] = 0.105171
$data[49,4] = 972697824
$data[49,5] = 46344605
$data[49,6] = 3.43311

$ws.Range("A2:G51").Value = $data
